# Update countries & provincias Spain
# - Kuwait's case count overtakes Marruecos / Argentina, so the ranked
#   rows 55-57 shift: Kuwait moves into row 55, Marruecos to row 56,
#   Argentina to row 57.
# - Refresh the daily numbers for several other ranked rows as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 (Emiratos Arabes Unidos) - updated daily figures
$ws.Range("B35").Value = 14163
$ws.Range("C35").Value = 564
$ws.Range("D35").Value = 2763
$ws.Range("E35").Value = 11274
$ws.Range("G35").Value = 7
$ws.Range("H35").Value = 126

# Row 36 (Polonia) - updated daily figures
$ws.Range("B36").Value = 13693
$ws.Range("C36").Value = 318
$ws.Range("D36").Value = 3945
$ws.Range("E36").Value = 9070
$ws.Range("G36").Value = 14
$ws.Range("H36").Value = 678

# Row 55 now becomes Kuwait (was Marruecos), with new figures
$ws.Range("A55").Value = "Kuwait"
$ws.Range("B55").Value = 4983
$ws.Range("C55").Value = 364
$ws.Range("D55").Value = 1776
$ws.Range("E55").Value = 3169
$ws.Range("F55").Value = 72
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 38

# Row 56 now becomes Marruecos (was Argentina), with new figures
$ws.Range("A56").Value = "Marruecos"
$ws.Range("B56").Value = 4880
$ws.Range("C56").Value = 151
$ws.Range("D56").Value = 1424
$ws.Range("E56").Value = 3282
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 174

# Row 57 now becomes Argentina (was Kuwait), with new figures
$ws.Range("A57").Value = "Argentina"
$ws.Range("B57").Value = 4681
$ws.Range("D57").Value = 1320
$ws.Range("E57").Value = 3124
$ws.Range("F57").Value = 157
$ws.Range("H57").Value = 237

# Row 76 (Bosnia y Herzegovina) - updated daily figures
$ws.Range("B76").Value = 1857
$ws.Range("C76").Value = 18
$ws.Range("D76").Value = 825
$ws.Range("E76").Value = 955
$ws.Range("G76").Value = 5
$ws.Range("H76").Value = 77

# Row 131 (Vietnam) - updated daily figures
$ws.Range("B131").Value = 271
$ws.Range("C131").Value = 1
$ws.Range("E131").Value = 52
